$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "319.89"
Set-TextValue "E2" "4.88%"
Set-TextValue "G2" "6"

Set-TextValue "D3" "36.04"
Set-TextValue "E3" "-0.15%"
Set-TextValue "G3" "6"

Set-TextValue "D4" "5.124"
Set-TextValue "E4" "0.88%"
Set-TextValue "G4" "6"

Set-TextValue "D5" "0.08191"
Set-TextValue "E5" "4.40%"
Set-TextValue "G5" "6"

Set-TextValue "D6" "2.147"
Set-TextValue "E6" "-2.69%"
Set-TextValue "G6" "6"

Set-TextValue "D7" "8.043"
Set-TextValue "E7" "1.58%"
Set-TextValue "G7" "6"

Set-TextValue "D8" "0.9253"
Set-TextValue "E8" "0.64%"
Set-TextValue "G8" "6"

Set-TextValue "E9" "4.80%"
Set-TextValue "G9" "6"

Set-TextValue "D10" "0.1892"
Set-TextValue "E10" "1.66%"
Set-TextValue "G10" "6"

Set-TextValue "D11" "0.09286"
Set-TextValue "E11" "7.13%"
Set-TextValue "G11" "6"

Set-TextValue "D12" "0.03594"
Set-TextValue "E12" "3.01%"
Set-TextValue "G12" "6"

Set-TextValue "D13" "0.09917"
Set-TextValue "E13" "-0.15%"
Set-TextValue "G13" "6"

Set-TextValue "D14" "0.001439"
Set-TextValue "E14" "0.78%"
Set-TextValue "G14" "6"

Set-TextValue "D15" "0.005690"
Set-TextValue "E15" "0.58%"
Set-TextValue "G15" "6"

Set-TextValue "D16" "3.460"
Set-TextValue "E16" "-0.09%"
Set-TextValue "G16" "6"

Set-TextValue "D17" "4.139"
Set-TextValue "E17" "0.94%"
Set-TextValue "G17" "6"

Set-TextValue "D18" "2.801"
Set-TextValue "E18" "13.98%"
Set-TextValue "G18" "6"

Set-TextValue "E19" "-1.52%"
Set-TextValue "G19" "6"

Set-TextValue "D20" "0.1330"
Set-TextValue "E20" "2.20%"
Set-TextValue "G20" "6"

Set-TextValue "D21" "5.099"
Set-TextValue "E21" "4.75%"
Set-TextValue "G21" "6"

Set-TextValue "D22" "0.2188"
Set-TextValue "E22" "-0.85%"
Set-TextValue "G22" "6"

Set-TextValue "D23" "0.04594"
Set-TextValue "E23" "0.71%"
Set-TextValue "G23" "6"

Set-TextValue "E24" "0.66%"
Set-TextValue "G24" "6"

Set-TextValue "D25" "0.004734"
Set-TextValue "E25" "-7.08%"
Set-TextValue "G25" "6"

Set-TextValue "D26" "0.0001299"
Set-TextValue "E26" "-7.39%"
Set-TextValue "G26" "6"

Set-TextValue "D27" "0.0004500"
Set-TextValue "E27" "-5.36%"
Set-TextValue "G27" "6"

Set-TextValue "G28" "6"

Set-TextValue "G29" "6"

Set-TextValue "G30" "6"

Set-TextValue "G31" "6"

Set-TextValue "G32" "6"

Set-TextValue "G33" "6"

Set-TextValue "G34" "6"

Set-TextValue "G35" "6"

Set-TextValue "G36" "6"

Set-TextValue "G37" "6"

Set-TextValue "G38" "6"

Set-TextValue "D39" "0.02011"
Set-TextValue "E39" "10.33%"
Set-TextValue "G39" "6"

Set-TextValue "D40" "0.04994"
Set-TextValue "E40" "4.89%"
Set-TextValue "G40" "6"

Set-TextValue "D41" "0.007828"
Set-TextValue "E41" "1.77%"
Set-TextValue "G41" "6"

Set-TextValue "E42" "0.27%"
Set-TextValue "G42" "6"

Set-TextValue "D43" "0.007810"
Set-TextValue "E43" "0.83%"
Set-TextValue "G43" "6"

Set-TextValue "D44" "0.002129"
Set-TextValue "E44" "-4.08%"
Set-TextValue "G44" "6"

Set-TextValue "D45" "0.01243"
Set-TextValue "E45" "11.05%"
Set-TextValue "G45" "6"

Set-TextValue "D46" "0.00006460"
Set-TextValue "E46" "2.86%"
Set-TextValue "G46" "6"

Set-TextValue "E47" "-0.19%"
Set-TextValue "G47" "6"

Set-TextValue "E48" "17.76%"
Set-TextValue "G48" "6"

Set-TextValue "D49" "0.001900"
Set-TextValue "E49" "-5.09%"
Set-TextValue "G49" "6"

Set-TextValue "D50" "0.00002100"
Set-TextValue "E50" "-0.19%"
Set-TextValue "G50" "6"

Set-TextValue "D51" "0.0002000"
Set-TextValue "E51" "-0.19%"
Set-TextValue "G51" "6"
